$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new tier-default rows (27-30): Adult Family Home / Childcare / Nursing Home / LTC
$ws.Range("A27").Value = "Adult Family Home"
$ws.Range("B27").Value = "afh"
$ws.Range("C27").Value = 4

$ws.Range("A28").Value = "Childcare"
$ws.Range("B28").Value = "child care"
$ws.Range("C28").Value = 3

$ws.Range("A29").Value = "Nursing Home"
$ws.Range("B29").Value = "ltcf"
$ws.Range("C29").Value = 1

$ws.Range("A30").Value = "LTC"
$ws.Range("B30").Value = "ltcf"
$ws.Range("C30").Value = 1

# Update the current selection to match the saved view state
$ws.Range("C34").Select()
